# registro.xlsx - "add - campo para fecha de inscripcion a materia"
#
# Summary of the change:
#  - Materias: add a new row (Materias-75e70e / LOG / Logica computacional)
#  - Inscripciones: replace the single existing record with two new
#    enrolment records (now living in rows 2-3, row 1 left blank) and add a
#    new "fecha de inscripcion" (enrolment date) column E, stored as text.
#  - Inscripciones becomes the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Materias: append the new subject row
# ---------------------------------------------------------------------
$wsMaterias = $wb.Worksheets.Item("Materias")
$wsMaterias.Cells.Item(4, 1).Value = "Materias-75e70e"
$wsMaterias.Cells.Item(4, 2).Value = "LOG"
$wsMaterias.Cells.Item(4, 3).Value = "Logica computacional"

# ---------------------------------------------------------------------
# Inscripciones: drop the old record, write the two new ones starting
# at row 2, and add the enrolment-date column (E)
# ---------------------------------------------------------------------
$wsInsc = $wb.Worksheets.Item("Inscripciones")

# the old row 1 record is removed entirely; the row is left empty
$wsInsc.Rows.Item(1).ClearContents()

# stamp the new enrolment-date column (E) as Text for the whole column,
# then scrub the leftover formatting from the (empty) row-1 cell so row 1
# stays completely blank
$wsInsc.Columns.Item(5).NumberFormat = "@"
$wsInsc.Cells.Item(1, 5).Style = "Normal"
$wsInsc.Cells.Item(1, 5).ClearContents()

# -- row 2 --
$wsInsc.Cells.Item(2, 1).Value = "Inscripciones-c2a095"
$wsInsc.Cells.Item(2, 2).Value = "Estudiantes-0e66c5"
$wsInsc.Cells.Item(2, 3).Value = "Materias-0dd25c"

$wsInsc.Cells.Item(2, 4).NumberFormat = "@"
$wsInsc.Cells.Item(2, 4).Value = "5.0"
$wsInsc.Cells.Item(2, 4).Style = "Normal"

$wsInsc.Cells.Item(2, 5).Value = "2025-09-11 09:58:54"
$wsInsc.Cells.Item(2, 5).Style = "Normal"

# -- row 3 --
$wsInsc.Cells.Item(3, 1).Value = "Inscripciones-77ee5a"
$wsInsc.Cells.Item(3, 2).Value = "Estudiantes-0e66c5"
$wsInsc.Cells.Item(3, 3).Value = "Materias-6ea3a0"

$wsInsc.Cells.Item(3, 4).NumberFormat = "@"
$wsInsc.Cells.Item(3, 4).Value = "7.0"
$wsInsc.Cells.Item(3, 4).Style = "Normal"

$wsInsc.Cells.Item(3, 5).Value = "2025-09-11 10:01:30"
$wsInsc.Cells.Item(3, 5).Style = "Normal"

# column widths for the Inscripciones sheet (COM ColumnWidth is in
# characters; the engine snaps to the same pixel grid real Excel uses)
$wsInsc.Columns.Item(1).ColumnWidth = 17.6328125 - (5 / 6)
$wsInsc.Columns.Item(2).ColumnWidth = 15 - (5 / 6)
$wsInsc.Columns.Item(3).ColumnWidth = 13.81640625 - (5 / 6)
$wsInsc.Columns.Item(5).ColumnWidth = 26.1796875 - (5 / 6)

# page orientation
$wsInsc.PageSetup.Orientation = 1

# make Inscripciones the active sheet/tab and select the (now empty)
# first row, matching the recorded UI state
$wsInsc.Activate()
$wsInsc.Rows.Item(1).Select() | Out-Null
